# Applies the cryptocurrency price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: a handful of "Price" (column D) values look like plain decimals (e.g. "0.540",
# "6.41", "8.15" ...). Assigning them straight to .Value would make Excel auto-convert
# them to numbers and drop the significant trailing zero / precision, which would not
# match the source data (these columns are plain text in the workbook). Prefixing the
# value with a leading apostrophe forces Excel to keep them as text, exactly like typing
# "'0.540" into the cell by hand.

$ws.Range("D2").Value = "69.156.41"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.758.25"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'602.48"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'167.31"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "3.758.10"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.540"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").Value = "'6.41"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'38.09"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "4.385.82"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "3.756.89"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "69.160.23"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "'17.34"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  +12.67%  "
$ws.Range("D22").Value = "'493.63"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  +7.00%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'12.32"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'2.48"
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").Value = "'31.63"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "3.903.26"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "3.689.57"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'6.00"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").Value = "'0.327"
$ws.Range("E42").Value = "  +5.47%  "
$ws.Range("D43").Value = "'430.30"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "'48.66"
$ws.Range("D45").Value = "'1.99"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D48").Value = "'40.40"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'141.31"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").Value = "2.794.49"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("E51").Value = "  +0.45%  "
